# Remove the CRP / Reynolds Risk Score rows (women + men) from the CVD
# prediction-models table. These occupy rows 4 and 5 (the SCORE row that
# follows shifts up to become the new row 4).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(4).Delete()
$ws.Rows.Item(4).Delete()

# Restore the view state recorded in the saved workbook: selection moves to
# F11 (and the forced top-left scroll anchor on A4 is cleared as a result).
$ws.Range("F11").Select()
